$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.9914139999999999
$ws.Range("H2").Value = 2.974242
$ws.Range("I2").Value = 0.5837683597777463
$ws.Range("J2").Value = 0.5837683597777463
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.5698856666666666
$ws.Range("N2").Value = 1.709657
$ws.Range("Q2").Value = 0.5649926283326666
$ws.Range("R2").Value = 5.084933654994
$ws.Range("S2").Value = 0.5837683597777463
$ws.Range("T2").Value = 0.5837683597777463

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 0.5478883333333334
$ws.Range("H3").Value = 1.643665
$ws.Range("I3").Value = 0.3226098014465836
$ws.Range("J3").Value = 0.3226098014465836
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.5698856666666666
$ws.Range("N3").Value = 1.709657
$ws.Range("Q3").Value = 0.3122337081005556
$ws.Range("R3").Value = 2.810103372905
$ws.Range("S3").Value = 0.3226098014465836
$ws.Range("T3").Value = 0.3226098014465836

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 0.158998
$ws.Range("H4").Value = 0.476994
$ws.Range("I4").Value = 0.09362183877567003
$ws.Range("J4").Value = 0.09362183877567003
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5698856666666666
$ws.Range("N4").Value = 1.709657
$ws.Range("Q4").Value = 0.09061068122866665
$ws.Range("R4").Value = 0.815496131058
$ws.Range("S4").Value = 0.09362183877567003
$ws.Range("T4").Value = 0.09362183877567003
